# Apply the trade update described in the commit:
#   Trade #107 (momentum strategy) closed via early_exit.
#   A new trade #136 (MarketMaking strategy) opened.
# This touches: Summary, Strategy Status, All Trades, momentum, MarketMaking

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.43   # Current Capital
$wsSummary.Range("B4").Value = 0.54      # Total P&L $
$wsSummary.Range("B5").Value = 0.1       # Total P&L %
$wsSummary.Range("B6").Value = 107       # Total Trades
$wsSummary.Range("B8").Value = 39        # Losing Trades
$wsSummary.Range("B9").Value = 48.6      # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - row 11 = momentum
# ---------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C11").Value = 99.27
$wsStatus.Range("D11").Value = 25
$wsStatus.Range("E11").Value = -0.73
$wsStatus.Range("F11").Value = -0.73
$wsStatus.Range("G11").Value = 28

# ---------------------------------------------------------------
# All Trades sheet - row 108 = trade #107 (momentum) closed early
# Columns: A#, B Date, C Time, D Strategy, E Side, F Entry Price,
#          G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#          L Exit Reason, M Duration (min), N Entry Slippage,
#          O Exit Slippage, P Confidence, Q Entry Reason
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Range("G108").Value = 0.92
$wsAll.Range("H108").Value = "CLOSED"
$wsAll.Range("I108").Value = -3.1579
$wsAll.Range("J108").Value = -0.03
$wsAll.Range("K108").Value = 99.27
$wsAll.Range("L108").Value = "early_exit"
$wsAll.Range("M108").Value = 0.12

# All Trades sheet - new row 137 = trade #136 (MarketMaking) opened
$wsAll.Range("A137").Value = 136
# Write the date as text (not an auto-converted date serial): force the
# cell to Text format before assigning, then restore the default style
# (copied from a neighboring plain-text date cell in the same column).
$wsAll.Range("B137").NumberFormat = "@"
$wsAll.Range("B137").Value = "2026-02-18"
$wsAll.Range("B136").Copy()
$wsAll.Range("B137").PasteSpecial(-4122)
$wsAll.Range("C137").Value = "00:28:16"
$wsAll.Range("D137").Value = "MarketMaking"
$wsAll.Range("E137").Value = "DOWN"
$wsAll.Range("F137").Value = 0.95
$wsAll.Range("H137").Value = "OPEN"
$wsAll.Range("I137").Value = 0
$wsAll.Range("J137").Value = 0
$wsAll.Range("K137").Value = 99.47967800952271
$wsAll.Range("M137").Value = 0
$wsAll.Range("N137").Value = 0
$wsAll.Range("O137").Value = 0
$wsAll.Range("P137").Value = 0.65
$wsAll.Range("Q137").Value = "Wide spread capture: 392 bps vs avg 315 bps"

# ---------------------------------------------------------------
# momentum sheet - row 26 = trade #107 closed early
# Columns: A#, B Date, C Time, D Strategy, E Side, F Entry Price,
#          G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#          L Entry Slippage, M Exit Slippage, N Confidence,
#          O Entry Reason, P Exit Reason, Q Duration (min)
# ---------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
$wsMomentum.Range("G26").Value = 0.92
$wsMomentum.Range("H26").Value = "CLOSED"
$wsMomentum.Range("I26").Value = -3.1579
$wsMomentum.Range("J26").Value = -0.03
$wsMomentum.Range("K26").Value = 99.27
$wsMomentum.Range("P26").Value = "early_exit"
$wsMomentum.Range("Q26").Value = 0.12

# ---------------------------------------------------------------
# MarketMaking sheet - new row 57 = trade #136 opened
# ---------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("A57").Value = 136
$wsMM.Range("B57").NumberFormat = "@"
$wsMM.Range("B57").Value = "2026-02-18"
$wsMM.Range("B56").Copy()
$wsMM.Range("B57").PasteSpecial(-4122)
$wsMM.Range("C57").Value = "00:28:16"
$wsMM.Range("D57").Value = "MarketMaking"
$wsMM.Range("E57").Value = "DOWN"
$wsMM.Range("F57").Value = 0.95
$wsMM.Range("H57").Value = "OPEN"
$wsMM.Range("I57").Value = 0
$wsMM.Range("J57").Value = 0
$wsMM.Range("K57").Value = 99.47967800952271
$wsMM.Range("L57").Value = 0
$wsMM.Range("M57").Value = 0
$wsMM.Range("N57").Value = 0.65
$wsMM.Range("O57").Value = "Wide spread capture: 392 bps vs avg 315 bps"
$wsMM.Range("Q57").Value = 0
